# Updated cryptos list on Fri Aug 16 13:45:33 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.369.99"
$ws.Range("E2").Value = "  -1.30%  "

$ws.Range("D3").Value = "2.591.66"
$ws.Range("E3").Value = "  -1.86%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'519.13"
$ws.Range("E5").Value = "  -1.49%  "

$ws.Range("D6").Value = "'141.39"
$ws.Range("E6").Value = "  -2.26%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("E8").Value = "  -0.88%  "

$ws.Range("D9").Value = "2.613.63"
$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("E10").Value = "  -2.64%  "

$ws.Range("E11").Value = "  -2.93%  "

$ws.Range("E12").Value = "  -1.43%  "

$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").Value = "3.047.65"
$ws.Range("E14").Value = "  -2.13%  "

$ws.Range("D15").Value = "58.335.36"
$ws.Range("E15").Value = "  -1.37%  "

$ws.Range("D16").Value = "'20.30"
$ws.Range("E16").Value = "  -3.34%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.593.96"
$ws.Range("E17").Value = "  -2.64%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000133"
$ws.Range("E18").Value = "  -2.73%  "

$ws.Range("D19").Value = "'337.82"
$ws.Range("E19").Value = "  -1.20%  "

$ws.Range("D20").Value = "'4.31"
$ws.Range("E20").Value = "  -3.15%  "

$ws.Range("D21").Value = "'10.21"

$ws.Range("D22").Value = "'6.40"
$ws.Range("E22").Value = "  +0.93%  "

$ws.Range("D24").Value = "'65.52"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "'0.167"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").Value = "'0.402"
$ws.Range("E26").Value = "  -3.92%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "2.698.98"
$ws.Range("E28").Value = "  -2.06%  "

$ws.Range("D29").Value = "'7.02"
$ws.Range("E29").Value = "  -2.99%  "

$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0738"
$ws.Range("E31").Value = "  -7.72%  "

$ws.Range("E32").Value = "  -6.06%  "

$ws.Range("E33").Value = "  -2.38%  "

$ws.Range("D34").Value = "'18.74"
$ws.Range("E34").Value = "  -1.12%  "

$ws.Range("D35").Value = "'149.30"
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("D36").Value = "'3.97"
$ws.Range("E36").Value = "  -5.85%  "

$ws.Range("E37").Value = "  -6.55%  "

$ws.Range("D38").Value = "'0.853"
$ws.Range("E38").Value = "  -2.40%  "

$ws.Range("D39").Value = "'36.25"
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("E40").Value = "  -1.90%  "

$ws.Range("D41").Value = "'0.831"
$ws.Range("E41").Value = "  -10.61%  "

$ws.Range("D42").Value = "'3.51"
$ws.Range("E42").Value = "  -3.91%  "

$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("D44").Value = "'273.72"
$ws.Range("E44").Value = "  +0.60%  "

$ws.Range("D45").Value = "'0.603"
$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("E46").Value = "  +0.26%  "

$ws.Range("E47").Value = "  -2.63%  "

$ws.Range("E48").Value = "  -3.63%  "

$ws.Range("D49").Value = "'18.58"
$ws.Range("E49").Value = "  -4.60%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.971.03"
$ws.Range("E50").Value = "  -3.92%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'4.64"
$ws.Range("E51").Value = "  -3.10%  "
